# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet with new
# market figures pulled from the latest CoinRanking snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value into a cell. Some refreshed "Price" values look
# like plain numbers (e.g. "8.30", "0.706"); Excel would otherwise auto-convert those
# to numbers (dropping trailing zeros), so pin the cell to Text first, then restore the
# default "Normal" style afterwards so no visible formatting/style change is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "63.796.05"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.089.95"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "538.58"
$ws.Range("E5").Value = "  -2.55%  "
Set-TextValue "D6" "136.68"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.082.99"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  +4.13%  "
Set-TextValue "D14" "34.93"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.591.09"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "63.795.64"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "3.094.82"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.19%  "
Set-TextValue "D20" "488.87"
$ws.Range("E20").Value = "  -3.51%  "
Set-TextValue "D21" "13.54"
$ws.Range("E21").Value = "  -0.37%  "
Set-TextValue "D22" "0.706"
$ws.Range("E22").Value = "  -0.21%  "
Set-TextValue "D23" "7.21"
$ws.Range("E23").Value = "  -0.83%  "
Set-TextValue "D24" "80.04"
$ws.Range("E24").Value = "  +2.56%  "
Set-TextValue "D25" "12.27"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.92%  "
Set-TextValue "D28" "8.30"
$ws.Range("E28").Value = "  +0.34%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.02%  "
Set-TextValue "D30" "26.31"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -0.20%  "
Set-TextValue "D33" "2.41"
$ws.Range("E33").Value = "  -5.26%  "
Set-TextValue "D34" "57.19"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +5.52%  "
Set-TextValue "D36" "502.73"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "3.296.72"
$ws.Range("E38").Value = "  +6.98%  "
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("E40").Value = "  +0.84%  "
Set-TextValue "D41" "0.118"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  -2.35%  "
Set-TextValue "D44" "0.259"
$ws.Range("E44").Value = "  +2.17%  "
Set-TextValue "D46" "2.12"
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("E47").Value = "  +6.51%  "
$ws.Range("E48").Value = "  +3.06%  "
Set-TextValue "D49" "122.32"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  +2.36%  "
Set-TextValue "D51" "2.34"
$ws.Range("E51").Value = "  -16.45%  "
